# Penalty Reward System rework (unfinished per commit message):
# - Remove several weekly / monthly rows whose PO quantities are being
#   reworked into a new penalty/reward calculation.
# - Update two of the remaining weekly figures and the corresponding
#   monthly-rollup figures to their new values.
#
# Sheet "Weekly Quantity" (sheet index 1): delete the rows for order weeks
# 2023-06-21 (45088.99...), 2023-07-05 (45102.99...), 2023-07-19 (45116.99...),
# 2023-08-02 (45130.99...), 2024-03-23 (45361.99...), 2024-03-30 (45368.99...)
# Then update the row now holding 45095.99... (was 230 -> 110) and the row
# now holding 45123.99... (was 170 -> 130).
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

# Delete rows bottom-up so earlier row numbers stay valid as we go.
$wsWeekly.Rows.Item(27).Delete()
$wsWeekly.Rows.Item(26).Delete()
$wsWeekly.Rows.Item(8).Delete()
$wsWeekly.Rows.Item(6).Delete()
$wsWeekly.Rows.Item(5).Delete()
$wsWeekly.Rows.Item(3).Delete()

# After the deletions, row 3 = 45095.99999999999, row 4 = 45123.99999999999.
$wsWeekly.Cells.Item(3, 2).Value = 110
$wsWeekly.Cells.Item(4, 2).Value = 130

# Sheet "Monthly Trend" (sheet index 2): delete the row for order month
# 2024-04-30 (45382.99...), then update the rows for 45107.99... (was 340 ->
# 110) and 45138.99... (was 390 -> 170).
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsMonthly.Rows.Item(11).Delete()

$wsMonthly.Cells.Item(3, 2).Value = 110
$wsMonthly.Cells.Item(4, 2).Value = 170
